$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$h = $ws.Hyperlinks.Item(1)
$h | Get-Member
